# working update for permissions and all other admin features
#
# Rebuilds the sample user-upload template:
#  - header row gains "role" and "department" columns (and "group" moves to
#    the end)
#  - the sample data row becomes a single "student_1" record with a
#    mailto: hyperlink on the email cell
#  - the old second sample data row (row 3) is cleared out but keeps its
#    formatting
#  - window/selection state is refreshed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear out the old second data row's contents first, before we start
#     writing the new shared strings, so no stale strings from "jane" /
#     "Group 2" / etc. linger in the workbook ------------------------------
$ws.Range("A3:E3").ClearContents()

# --- header row: insert "role"/"department", keep "group" as the last
#     column --------------------------------------------------------------
$ws.Range("E1").Value = "role"
$ws.Range("F1").Value = "department"
$ws.Range("G1").Value = "group"

# copy the header style (bold, centered, wrapped) onto the two new header
# cells
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# --- data row 2: replace the sample "john" row with a "student_1" record -
$ws.Range("A2").Value = "student_1"
$ws.Range("B2").Value = "student_1"
$ws.Range("C2").Value = "stu1@gmail.com"
$ws.Range("D2").Value = "password"
$ws.Range("E2").Value = "student"
$ws.Range("F2").Value = "Computer Science"

# copy the data row style onto the new department cell
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# --- turn the e-mail cell into a real mailto: hyperlink -------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:stu1@gmail.com")

# --- row heights: header row now wraps to two lines like the data rows;
#     the now-empty third row goes back to the default height ------------
$ws.Rows.Item(1).RowHeight = 28.5
$ws.Rows.Item(3).AutoFit()

# --- window / selection state ---------------------------------------------
$win = $excel.Windows.Item(1)
$win.Left = 1747
$win.Top = 1747
$win.Width = 14401
$win.Height = 8183

$ws.Range("H2").Select()
